$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing D7 cell text (was "13:30-zeit", becomes "13:30-20:00")
$ws.Range("D7").Value = "13:30-20:00"

# Remember the time-format (h:mm) used by the "Additional Notes" column so new
# rows can reuse the exact same style as the existing ones.
$timeFormat = $ws.Range("D7").NumberFormat

# --- New log entries: 12.12.2023, 13.12.2023, 18.12.2023 ---

# Row 8 date first
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "13.12.2023"
$ws.Range("A8").ClearFormats()

# Row 7 work description
$ws.Range("C7").Value = "Continued on Server"

# Row 7 date
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "12.12.2023"
$ws.Range("A7").ClearFormats()

# Row 8 remaining cells
$ws.Range("B7").Value = 390
$ws.Range("B8").Value = 90
$ws.Range("C8").Value = "Continued on Server"
$ws.Range("D8").Value = "18:30-20:00"
$ws.Range("D8").NumberFormat = $timeFormat

# Row 9 (18.12.2023 entry)
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "18.12.2023"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = 210
$ws.Range("C9").Value = "Started on DB"
$ws.Range("D9").Value = "15:30-19:00"
$ws.Range("D9").NumberFormat = $timeFormat

# Update selection to match the target view state
$ws.Range("B8").Select()
